$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 45458364
$ws.Range("I86").Value = 2191.182
$ws.Range("J86").Value = 90914536
$ws.Range("K86").Value = 2191.182
$ws.Range("L86").Value = 90914536
$ws.Range("M86").Value = -1068.182
$ws.Range("N86").Value = -90916782
$ws.Range("H89").Value = 45458364
$ws.Range("I89").Value = 2191.182
$ws.Range("J89").Value = 90914536
$ws.Range("K89").Value = 10955.91
$ws.Range("L89").Value = 454572680
$ws.Range("M89").Value = -5339.91
$ws.Range("N89").Value = -454583912
$ws.Range("H113").Value = 2205.3713
$ws.Range("I113").Value = 2421.9443
$ws.Range("J113").Value = 1976.0588
$ws.Range("K113").Value = 2421.9443
$ws.Range("L113").Value = 1976.0588
$ws.Range("M113").Value = 832.0556999999999
$ws.Range("N113").Value = -8484.058800000001
$ws.Range("H121").Value = 704.3103599999999
$ws.Range("J121").Value = 697.2222
$ws.Range("L121").Value = 2091.6666
$ws.Range("N121").Value = -5585.6666
$ws.Range("H132").Value = 1524.081
$ws.Range("I132").Value = 1486.6129
$ws.Range("J132").Value = 1717.6666
$ws.Range("K132").Value = 4459.8387
$ws.Range("L132").Value = 5152.9998
$ws.Range("M132").Value = -1929.8387
$ws.Range("N132").Value = -10212.9998
$ws.Range("H133").Value = 35000
$ws.Range("J133").Value = 35000
$ws.Range("L133").Value = 35000
$ws.Range("N133").Value = -45120
$ws.Range("H137").Value = 7576950.5
$ws.Range("I137").Value = 1201.12
$ws.Range("J137").Value = 31251166
$ws.Range("K137").Value = 3603.36
$ws.Range("L137").Value = 93753498
$ws.Range("M137").Value = -1053.36
$ws.Range("N137").Value = -93758598
$ws.Range("H141").Value = 1018.6
$ws.Range("I141").Value = 909.55554
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 2728.66662
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 2451.33338
$ws.Range("N141").Value = -16360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16854.082
$ws.Range("I32").Value = 14823.695
$ws.Range("J32").Value = 54213.2
$ws.Range("K32").Value = 14823.695
$ws.Range("L32").Value = 54213.2
$ws.Range("M32").Value = -14536.695
$ws.Range("N32").Value = -54787.2
$ws.Range("H45").Value = 6710
$ws.Range("I45").Value = 7763.5
$ws.Range("J45").Value = 4302
$ws.Range("K45").Value = 7763.5
$ws.Range("L45").Value = 4302
$ws.Range("M45").Value = -7386.5
$ws.Range("N45").Value = -5056
$ws.Range("H61").Value = 1682.6923
$ws.Range("I61").Value = 1431.1428
$ws.Range("J61").Value = 1976.1666
$ws.Range("K61").Value = 1431.1428
$ws.Range("L61").Value = 1976.1666
$ws.Range("M61").Value = -1219.1428
$ws.Range("N61").Value = -2400.1666
$ws.Range("H63").Value = 2332.8572
$ws.Range("I63").Value = 2332.8572
$ws.Range("K63").Value = 2332.8572
$ws.Range("M63").Value = -1646.8572
$ws.Range("H66").Value = 2332.8572
$ws.Range("I66").Value = 2332.8572
$ws.Range("K66").Value = 11664.286
$ws.Range("M66").Value = -8232.286
$ws.Range("H74").Value = 1057.2572
$ws.Range("I74").Value = 1268.8422
$ws.Range("J74").Value = 806
$ws.Range("K74").Value = 1268.8422
$ws.Range("L74").Value = 806
$ws.Range("M74").Value = -394.8422
$ws.Range("N74").Value = -2554
$ws.Range("H77").Value = 1057.2572
$ws.Range("I77").Value = 1268.8422
$ws.Range("J77").Value = 806
$ws.Range("K77").Value = 6344.211
$ws.Range("L77").Value = 4030
$ws.Range("M77").Value = -1976.211
$ws.Range("N77").Value = -12766
$ws.Range("H136").Value = 1682.6923
$ws.Range("I136").Value = 1431.1428
$ws.Range("J136").Value = 1976.1666
$ws.Range("K136").Value = 4293.428400000001
$ws.Range("L136").Value = 5928.4998
$ws.Range("M136").Value = -1743.428400000001
$ws.Range("N136").Value = -11028.4998
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2322.4443
$ws.Range("J86").Value = 2418.4546
$ws.Range("L86").Value = 2418.4546
$ws.Range("N86").Value = -4664.4546
$ws.Range("H89").Value = 2322.4443
$ws.Range("J89").Value = 2418.4546
$ws.Range("L89").Value = 12092.273
$ws.Range("N89").Value = -23324.273
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1709.3823
$ws.Range("I31").Value = 1096.6
$ws.Range("J31").Value = 2193.158
$ws.Range("K31").Value = 1096.6
$ws.Range("L31").Value = 2193.158
$ws.Range("M31").Value = -801.5999999999999
$ws.Range("N31").Value = -2783.158
$ws.Range("H34").Value = 1709.3823
$ws.Range("I34").Value = 1096.6
$ws.Range("J34").Value = 2193.158
$ws.Range("K34").Value = 1096.6
$ws.Range("L34").Value = 2193.158
$ws.Range("M34").Value = -894.5999999999999
$ws.Range("N34").Value = -2597.158
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 521.53845
$ws.Range("I15").Value = 182.85715
$ws.Range("J15").Value = 916.6667
$ws.Range("K15").Value = 548.5714499999999
$ws.Range("L15").Value = 2750.0001
$ws.Range("M15").Value = -408.5714499999999
$ws.Range("N15").Value = -3030.0001
$ws.Range("H17").Value = 565
$ws.Range("I17").Value = 350.5
$ws.Range("J17").Value = 672.25
$ws.Range("K17").Value = 1051.5
$ws.Range("L17").Value = 2016.75
$ws.Range("M17").Value = -882.5
$ws.Range("N17").Value = -2354.75
$ws.Range("H68").Value = 1342.0385
$ws.Range("I68").Value = 765.0333000000001
$ws.Range("J68").Value = 1702.6666
$ws.Range("K68").Value = 2295.0999
$ws.Range("L68").Value = 5107.9998
$ws.Range("M68").Value = -1484.0999
$ws.Range("N68").Value = -6729.9998
$ws.Range("H71").Value = 1342.0385
$ws.Range("I71").Value = 765.0333000000001
$ws.Range("J71").Value = 1702.6666
$ws.Range("K71").Value = 6885.2997
$ws.Range("L71").Value = 15323.9994
$ws.Range("M71").Value = -2829.2997
$ws.Range("N71").Value = -23435.9994
$ws.Range("H107").Value = 1764.3334
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1764.3334
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 5293.0002
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -9133.0002
$ws.Range("H113").Value = 161840.42
$ws.Range("I113").Value = 421.42856
$ws.Range("J113").Value = 182384.66
$ws.Range("K113").Value = 1264.28568
$ws.Range("L113").Value = 547153.98
$ws.Range("M113").Value = 905.71432
$ws.Range("N113").Value = -551493.98
$ws.Range("H117").Value = 2864.3845
$ws.Range("I117").Value = 962.3333
$ws.Range("J117").Value = 3435
$ws.Range("K117").Value = 2886.9999
$ws.Range("L117").Value = 10305
$ws.Range("M117").Value = 555.0001000000002
$ws.Range("N117").Value = -17189
$ws.Range("H131").Value = 30255.861
$ws.Range("I131").Value = 101515
$ws.Range("J131").Value = 2848.5
$ws.Range("K131").Value = 304545
$ws.Range("L131").Value = 8545.5
$ws.Range("M131").Value = -299505
$ws.Range("N131").Value = -18625.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2533.4443
$ws.Range("I97").Value = 2827.1428
$ws.Range("J97").Value = 1505.5
$ws.Range("K97").Value = 2827.1428
$ws.Range("L97").Value = 1505.5
$ws.Range("M97").Value = -2331.1428
$ws.Range("N97").Value = -2497.5
$ws.Range("H128").Value = 38211.43
$ws.Range("J128").Value = 38211.43
$ws.Range("L128").Value = 38211.43
$ws.Range("N128").Value = -48171.43
$ws.Range("H130").Value = 54980
$ws.Range("J130").Value = 54980
$ws.Range("L130").Value = 54980
$ws.Range("N130").Value = -65020
$ws.Range("H133").Value = 41669.23
$ws.Range("J133").Value = 41669.23
$ws.Range("L133").Value = 41669.23
$ws.Range("N133").Value = -51789.23
$ws.Range("H135").Value = 39911.11
$ws.Range("J135").Value = 39911.11
$ws.Range("L135").Value = 39911.11
$ws.Range("N135").Value = -50051.11
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10740334
$ws.Range("I40").Value = 11086667
$ws.Range("K40").Value = 11086667
$ws.Range("M40").Value = -11086531
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 52561.12
$ws.Range("J135").Value = 52561.12
$ws.Range("L135").Value = 52561.12
$ws.Range("N135").Value = -62701.12
$ws.Range("H136").Value = 2530.8774
$ws.Range("I136").Value = 3005.111
$ws.Range("J136").Value = 1948.8636
$ws.Range("K136").Value = 9015.332999999999
$ws.Range("L136").Value = 5846.5908
$ws.Range("M136").Value = -6465.332999999999
$ws.Range("N136").Value = -10946.5908
